$wb = $excel.ActiveWorkbook

$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodes = $wb.Worksheets.Item("Codes")
$wsExtensionSchemes = $wb.Worksheets.Item("ExtensionSchemes")
$wsExtensionsTest = $wb.Worksheets.Item("Extensions_test")

# Rename sheets to the new, dynamic (codevalue-suffixed) names.
$wsCodes.Name = "Codes_exttest1"
$wsExtensionSchemes.Name = "ExtensionSchemes_exttest1"

# New header/value cells on the CodeSchemes sheet (columns K & L), and on the
# ExtensionSchemes sheet (column H), referencing the sheet names used for the
# Codes / ExtensionSchemes / Extensions sheets of this codescheme.
# Order matches the order the new shared strings were first written.
$wsCodeSchemes.Range("L1").Value = "EXTENSIONSCHEMESSHEET"
$wsCodeSchemes.Range("K1").Value = "CODESSHEET"
$wsExtensionSchemes.Range("H1").Value = "EXTENSIONSSHEET"
$wsCodeSchemes.Range("K2").Value = "Codes_exttest1"
$wsCodeSchemes.Range("L2").Value = "ExtensionSchemes_exttest1"
$wsExtensionSchemes.Range("H2").Value = "Extensions_test"

# Match the style (text number format) used by the other "sheet name" style cells.
$wsCodeSchemes.Range("K2").NumberFormat = "@"
$wsCodeSchemes.Range("L2").NumberFormat = "@"

# New column width for the new CodeSchemes column K.
$wsCodeSchemes.Columns.Item(11).ColumnWidth = 12.83

# Move the active tab / selections: ExtensionSchemes_exttest1 becomes the
# tabSelected sheet (was Codes before).
$wsCodeSchemes.Range("L3").Select() | Out-Null
$wsExtensionSchemes.Range("H2").Select() | Out-Null
$wsExtensionSchemes.Activate()
